$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "n" row values for the X0 and X1 columns look purely numeric, so force
# them to be stored as text (matching the original shared-string type)
# by marking the cells as Text before assigning.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "     72"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "     78"

# Overall column (B)
$ws.Range("B4").Value = "      3 ( 2.0) "
$ws.Range("B5").Value = "     16 (10.7) "
$ws.Range("B6").Value = "     38 (25.3) "
$ws.Range("B7").Value = "     41 (27.3) "
$ws.Range("B8").Value = "     52 (34.7) "
$ws.Range("B9").Value = "   0.82 (0.53)"
$ws.Range("B10").Value = "4202.76 (4467.12)"

# X0..not.D.E.F column (C)
$ws.Range("C4").Value = "      2 ( 2.8) "
$ws.Range("C5").Value = "      5 ( 6.9) "
$ws.Range("C6").Value = "     18 (25.0) "
$ws.Range("C7").Value = "     23 (31.9) "
$ws.Range("C8").Value = "     24 (33.3) "
$ws.Range("C9").Value = "   0.95 (0.57)"
$ws.Range("C10").Value = "4825.92 (4825.38)"

# X1..best.D.E.F column (D)
$ws.Range("D4").Value = "      1 ( 1.3) "
$ws.Range("D5").Value = "     11 (14.1) "
$ws.Range("D6").Value = "     20 (25.6) "
$ws.Range("D7").Value = "     18 (23.1) "
$ws.Range("D8").Value = "     28 (35.9) "
$ws.Range("D9").Value = "   0.71 (0.46)"
$ws.Range("D10").Value = "3627.54 (4055.41)"
